$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Cells.Item(2, 4)
$cell.NumberFormat = "@"
$cell.Value = '25.776.48'
$cell.Style = "Normal"
$ws.Cells.Item(2, 5).Value = '  -0.18%  '

$cell = $ws.Cells.Item(3, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.634.99'
$cell.Style = "Normal"
$ws.Cells.Item(3, 5).Value = '  +0.06%  '

$ws.Cells.Item(4, 5).Value = '  -0.18%  '

$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = '215.47'
$cell.Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  +0.22%  '

$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.501'
$cell.Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  -0.78%  '

$ws.Cells.Item(7, 5).Value = '  -0.11%  '

$ws.Cells.Item(8, 5).Value = '  -0.04%  '

$cell = $ws.Cells.Item(9, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.0634'
$cell.Style = "Normal"
$ws.Cells.Item(9, 5).Value = '  -1.13%  '

$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = "@"
$cell.Value = '19.56'
$cell.Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  -1.57%  '

$cell = $ws.Cells.Item(11, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.0791'
$cell.Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  +1.36%  '

$ws.Cells.Item(12, 5).Value = '  +0.18%  '

$cell = $ws.Cells.Item(13, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.861.29'
$cell.Style = "Normal"
$ws.Cells.Item(13, 5).Value = '  +0.13%  '

$cell = $ws.Cells.Item(14, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.639.61'
$cell.Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  -0.21%  '

$cell = $ws.Cells.Item(15, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.562'
$cell.Style = "Normal"
$ws.Cells.Item(15, 5).Value = '  +0.83%  '

$ws.Cells.Item(16, 5).Value = '  -0.57%  '

$cell = $ws.Cells.Item(17, 4)
$cell.NumberFormat = "@"
$cell.Value = '63.21'
$cell.Style = "Normal"
$ws.Cells.Item(17, 5).Value = '  +0.32%  '

$cell = $ws.Cells.Item(18, 4)
$cell.NumberFormat = "@"
$cell.Value = '25.814.47'
$cell.Style = "Normal"
$ws.Cells.Item(18, 5).Value = '  -0.01%  '

$ws.Cells.Item(19, 5).Value = '  -0.17%  '

$cell = $ws.Cells.Item(20, 4)
$cell.NumberFormat = "@"
$cell.Value = '4.46'
$cell.Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  +1.86%  '

$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = "@"
$cell.Value = '192.33'
$cell.Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  -0.85%  '

$ws.Cells.Item(22, 5).Value = '  +0.65%  '

$ws.Cells.Item(23, 5).Value = '  +2.14%  '

$ws.Cells.Item(24, 5).Value = '  +4.19%  '

$ws.Cells.Item(25, 5).Value = '  -0.10%  '

$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = "@"
$cell.Value = '141.63'
$cell.Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  +1.59%  '

$ws.Cells.Item(27, 5).Value = '  +1.10%  '

$cell = $ws.Cells.Item(28, 4)
$cell.NumberFormat = "@"
$cell.Value = '6.90'
$cell.Style = "Normal"
$ws.Cells.Item(28, 5).Value = '  +1.16%  '

$ws.Cells.Item(29, 5).Value = '  -0.43%  '

$ws.Cells.Item(30, 5).Value = '  -0.04%  '

$ws.Cells.Item(31, 5).Value = '  -0.39%  '

$ws.Cells.Item(32, 5).Value = '  +0.19%  '

$ws.Cells.Item(33, 5).Value = '  -0.39%  '

$ws.Cells.Item(34, 5).Value = '  -0.43%  '

$ws.Cells.Item(35, 5).Value = '  -0.46%  '

$cell = $ws.Cells.Item(36, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.903'
$cell.Style = "Normal"
$ws.Cells.Item(36, 5).Value = '  +0.22%  '

$cell = $ws.Cells.Item(37, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.132.15'
$cell.Style = "Normal"
$ws.Cells.Item(37, 5).Value = '  +1.37%  '

$cell = $ws.Cells.Item(38, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.52'
$cell.Style = "Normal"
$ws.Cells.Item(38, 5).Value = '  -1.94%  '

$cell = $ws.Cells.Item(39, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.544'
$cell.Style = "Normal"
$ws.Cells.Item(39, 5).Value = '  -0.87%  '

$ws.Cells.Item(40, 5).Value = '  -0.50%  '

$ws.Cells.Item(41, 5).Value = '  +0.11%  '

$ws.Cells.Item(42, 5).Value = '  +0.69%  '

$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = "@"
$cell.Value = '5.57'
$cell.Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  +0.70%  '

$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = "@"
$cell.Value = '100.50'
$cell.Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  +1.31%  '

$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.805'
$cell.Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  +0.78%  '

$cell = $ws.Cells.Item(46, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.770.08'
$cell.Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  -0.04%  '

$ws.Cells.Item(47, 5).Value = '  +3.86%  '

$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = "@"
$cell.Value = '55.32'
$cell.Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  -0.30%  '

$ws.Cells.Item(49, 5).Value = '  -0.92%  '

$ws.Cells.Item(50, 5).Value = '  -0.24%  '

$ws.Cells.Item(51, 5).Value = '  +4.33%  '
